# Combines Whole #'s For Everyone Properly
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the first four staff rows (3-6) with new names/values ---
$ws.Range("A3").Value = "Maggie  Farrell"
$ws.Range("B3").Value = 38
$ws.Range("C3").Value = 20
$ws.Range("D3").Value = 4

$ws.Range("A4").Value = "Makayla Baca"
$ws.Range("B4").Value = 24
$ws.Range("C4").Value = 7
$ws.Range("D4").Value = 3

$ws.Range("A5").Value = "Justyne Martinez "
$ws.Range("B5").Value = 25
$ws.Range("C5").Value = 15
$ws.Range("D5").Value = 4

$ws.Range("A6").Value = "Vy Torino"
$ws.Range("B6").Value = 17
$ws.Range("C6").Value = 6
$ws.Range("D6").Value = 6

# --- Append the combined (deduplicated, alphabetised) staff summary table ---
# Row 19 is intentionally left blank (gap before the new table).

$ws.Range("A20").Value = "Aminah Avalos"
$ws.Range("B20").Value = 34
$ws.Range("C20").Value = 14
$ws.Range("D20").Value = 8

$ws.Range("A21").Value = "Chrissy Cummings"
$ws.Range("B21").Value = 29
$ws.Range("C21").Value = 24
$ws.Range("D21").Value = 6

$ws.Range("A22").Value = "Danielle Mai"
$ws.Range("B22").Value = 9
$ws.Range("C22").Value = 7
$ws.Range("D22").Value = 1

$ws.Range("A23").Value = "Izzy Kruis"
$ws.Range("B23").Value = 44
$ws.Range("C23").Value = 24
$ws.Range("D23").Value = 6

$ws.Range("A24").Value = "Jasmine Saiz"
$ws.Range("B24").Value = 56
$ws.Range("C24").Value = 27
$ws.Range("D24").Value = 10

$ws.Range("A25").Value = "Justyne Martinez "
$ws.Range("B25").Value = 25
$ws.Range("C25").Value = 15
$ws.Range("D25").Value = 4

$ws.Range("A26").Value = "Karen Trevizo"
$ws.Range("B26").Value = 27
$ws.Range("C26").Value = 21
$ws.Range("D26").Value = 1

$ws.Range("A27").Value = "Krisdee Martinez"
$ws.Range("B27").Value = 39
$ws.Range("C27").Value = 19
$ws.Range("D27").Value = 10

$ws.Range("A28").Value = "Maggie  Farrell"
$ws.Range("B28").Value = 38
$ws.Range("C28").Value = 20
$ws.Range("D28").Value = 4

$ws.Range("A29").Value = "Makayla Baca"
$ws.Range("B29").Value = 24
$ws.Range("C29").Value = 7
$ws.Range("D29").Value = 3

$ws.Range("A30").Value = "Matthew Young"
$ws.Range("B30").Value = 1
$ws.Range("C30").Value = 1
$ws.Range("D30").Value = 0

$ws.Range("A31").Value = "Vy Torino"
$ws.Range("B31").Value = 17
$ws.Range("C31").Value = 6
$ws.Range("D31").Value = 6
